$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at row 271 (shifts existing rows 271+ down by 2)
$ws.Rows.Item(271).Insert()
$ws.Rows.Item(271).Insert()

# Fill new row 271
$ws.Cells.Item(271,1).Value2 = 8
$ws.Cells.Item(271,2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(271,3).Value2 = "Coquimbo"
$ws.Cells.Item(271,4).Value2 = 44489
$ws.Cells.Item(271,5).Value2 = 4
$ws.Cells.Item(271,6).Value2 = 100112045
$ws.Cells.Item(271,7).Value2 = "Zapallo"
$ws.Cells.Item(271,8).Value2 = "Camote"
$ws.Cells.Item(271,9).Value2 = "1a nueva(o)"
$ws.Cells.Item(271,10).Value2 = 800
$ws.Cells.Item(271,11).Value2 = 700
$ws.Cells.Item(271,12).Value2 = 750
$ws.Cells.Item(271,13).Value2 = 725
$ws.Cells.Item(271,14).Value2 = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(271,15).Value2 = "Perú"
$ws.Cells.Item(271,16).Value2 = 725
$ws.Cells.Item(271,17).Value2 = 1
$ws.Cells.Item(271,18).Value2 = "Hortaliza"

# Fill new row 272
$ws.Cells.Item(272,1).Value2 = 8
$ws.Cells.Item(272,2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(272,3).Value2 = "Coquimbo"
$ws.Cells.Item(272,4).Value2 = 44489
$ws.Cells.Item(272,5).Value2 = 4
$ws.Cells.Item(272,6).Value2 = 100112045
$ws.Cells.Item(272,7).Value2 = "Zapallo"
$ws.Cells.Item(272,8).Value2 = "Camote"
$ws.Cells.Item(272,9).Value2 = "2a nueva(o)"
$ws.Cells.Item(272,10).Value2 = 500
$ws.Cells.Item(272,11).Value2 = 600
$ws.Cells.Item(272,12).Value2 = 650
$ws.Cells.Item(272,13).Value2 = 625
$ws.Cells.Item(272,14).Value2 = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(272,15).Value2 = "Perú"
$ws.Cells.Item(272,16).Value2 = 625
$ws.Cells.Item(272,17).Value2 = 1
$ws.Cells.Item(272,18).Value2 = "Hortaliza"

Write-Output "done"